# "array de teste no lab" - adds a new "Validação" worksheet that mirrors the
# existing "Demonstração" sheet's layout/header, then drops in a fresh pair of
# test-rig sensor readings (250-50 / 168-45) to validate against.

$wb = $excel.ActiveWorkbook
$demo = $wb.Worksheets("Demonstração")

# New sheet goes right after "Demonstração", at the end of the tab strip.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $demo)
$ws.Name = "Validação"

# Carry over the header row + formatting (fonts, alignment, widths inherited
# from the styles already used on "Demonstração") instead of rebuilding it by
# hand - this keeps every style index identical to the sibling sheets.
$demo.Range("A1:M3").Copy($ws.Range("A1"))

# Row 2 - new "250-50" Deformação sensor (channel 2) read in the lab.
$ws.Range("A2").Value = "250-50"
$ws.Range("B2").Value = 1550.4313083333336
$ws.Range("D2").Value = 2
$ws.Range("K2").Value = 20.409904278669277

# Row 3 - new "168-45" Temperatura sensor (channel 2).
$ws.Range("A3").Value = "168-45"
$ws.Range("B3").Value = 1545.2719999999999
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = 34.6
$ws.Range("G3").Value = -0.7

# I3 is left blank but picks up a distinct (underlined, centered, auto-color)
# style - stub formatting for a value to be filled in later.
$ws.Range("I3").Font.Underline = 2
$ws.Range("I3").HorizontalAlignment = -4108

# "Demonstração" had row 1 selected while reviewing headers alongside the new
# sheet; its tab is no longer the active one once we flip to "Validação".
[void]$demo.Activate()
[void]$demo.Rows("1:1").Select()

# Leave the new sheet active/selected, cursor parked at L5.
[void]$ws.Activate()
[void]$ws.Range("L5").Select()
